$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("constants")
$ws1.Range("B15").Value = 26.24
